# Applies the "Updated cryptos list" data refresh to Sheet1 (columns B-E, rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain cell writes. A handful of Price (column D) values are bare decimal numbers
# ("211.97", "0.487", "1.00", ...) - same as in the source workbook, those are stored
# as text, not numbers. Range.Value auto-converts a bare numeric-looking string to a
# Double, which would both change the cell's type and mangle values like "1.00" -> 1 or
# introduce float noise ("0.487" -> 0.48699999999999999). Prefixing with a leading
# apostrophe is the standard Excel "force text" marker (same as a user typing '0.487 into
# a cell) - it keeps the stored value as the exact text "0.487" with no format change.

$ws.Range("D2").Value = '25.941.44'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '1.616.75'
$ws.Range("E3").Value = '  -1.06%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''211.97'
$ws.Range("E5").Value = '  -1.04%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = '''0.487'
$ws.Range("E7").Value = '  -3.57%  '
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("E9").Value = '  -1.60%  '
$ws.Range("D10").Value = '''18.23'
$ws.Range("E10").Value = '  -2.05%  '
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("D12").Value = '1.842.17'
$ws.Range("E12").Value = '  -1.02%  '
$ws.Range("D13").Value = '1.614.83'
$ws.Range("E13").Value = '  -1.44%  '
$ws.Range("D14").Value = '''4.12'
$ws.Range("E14").Value = '  -2.16%  '
$ws.Range("D15").Value = '''0.519'
$ws.Range("E15").Value = '  -2.20%  '
$ws.Range("D16").Value = '25.949.61'
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("D17").Value = '''61.63'
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("D18").Value = '0.0₃0733'
$ws.Range("E18").Value = '  -1.71%  '
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = '''191.35'
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("D21").Value = '''4.23'
$ws.Range("E21").Value = '  -0.86%  '
$ws.Range("D22").Value = '''9.47'
$ws.Range("E22").Value = '  -1.26%  '
$ws.Range("D23").Value = '''6.01'
$ws.Range("E23").Value = '  -2.20%  '
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("D25").Value = '''143.61'
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -3.26%  '
$ws.Range("D28").Value = '''6.61'
$ws.Range("E28").Value = '  -2.25%  '
$ws.Range("D29").Value = '''15.17'
$ws.Range("E29").Value = '  -0.63%  '
$ws.Range("D30").Value = '''1.22'
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("D31").Value = '''0.0475'
$ws.Range("E31").Value = '  -2.05%  '
$ws.Range("E32").Value = '  -1.50%  '
$ws.Range("E33").Value = '  -2.91%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '''2.40'
$ws.Range("E34").Value = '  -1.29%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").Value = '''1.48'
$ws.Range("E35").Value = '  -1.27%  '
$ws.Range("D36").Value = '1.127.25'
$ws.Range("E36").Value = '  -0.29%  '
$ws.Range("D37").Value = '''0.824'
$ws.Range("E37").Value = '  -5.95%  '
$ws.Range("E38").Value = '  -1.65%  '
$ws.Range("D39").Value = '''0.517'
$ws.Range("E39").Value = '  -1.78%  '
$ws.Range("E40").Value = '  -1.50%  '
$ws.Range("D41").Value = '''97.41'
$ws.Range("E41").Value = '  -1.48%  '
$ws.Range("D42").Value = '1.754.08'
$ws.Range("E42").Value = '  -0.87%  '
$ws.Range("E43").Value = '  -4.15%  '
$ws.Range("D44").Value = '''5.09'
$ws.Range("E44").Value = '  -4.03%  '
$ws.Range("E45").Value = '  -0.23%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '''53.84'
$ws.Range("E46").Value = '  -2.73%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '''1.48'
$ws.Range("E47").Value = '  -0.53%  '
$ws.Range("E48").Value = '  -2.27%  '
$ws.Range("D49").Value = '''0.412'
$ws.Range("E49").Value = '  -0.60%  '
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").Value = '''1.00'
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''7.42'
$ws.Range("E51").Value = '  -1.92%  '
